$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 417-418 (pushing the existing 417:441 block down to 419:443),
# then populate them with a new weekly record (date 44516) for "Apio" at
# Mercado Mayorista Lo Valledor de Santiago, following the same row layout
# as the rest of the dataset.
$ws.Rows("417:418").Insert()

# Row 417 - Primera
$ws.Range("A417").Value = 6
$ws.Range("B417").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C417").Value = "Metropolitana"
$ws.Range("D417").Value = 44516
$ws.Range("E417").Value = 13
$ws.Range("F417").Value = 100112017
$ws.Range("G417").Value = "Apio"
$ws.Range("H417").Value = "Americana (o)"
$ws.Range("I417").Value = "Primera"
$ws.Range("J417").Value = 2200
$ws.Range("K417").Value = 5000
$ws.Range("L417").Value = 6000
$ws.Range("M417").Value = 5364
$ws.Range("N417").Value = "`$/docena de matas"
$ws.Range("O417").Value = "Región de Coquimbo"
$ws.Range("P417").Value = 894
$ws.Range("Q417").Value = 6
$ws.Range("R417").Value = "Hortaliza"

# Row 418 - Segunda
$ws.Range("A418").Value = 6
$ws.Range("B418").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C418").Value = "Metropolitana"
$ws.Range("D418").Value = 44516
$ws.Range("E418").Value = 13
$ws.Range("F418").Value = 100112017
$ws.Range("G418").Value = "Apio"
$ws.Range("H418").Value = "Americana (o)"
$ws.Range("I418").Value = "Segunda"
$ws.Range("J418").Value = 800
$ws.Range("K418").Value = 4000
$ws.Range("L418").Value = 4000
$ws.Range("M418").Value = 4000
$ws.Range("N418").Value = "`$/docena de matas"
$ws.Range("O418").Value = "Región de Coquimbo"
$ws.Range("P418").Value = 667
$ws.Range("Q418").Value = 6
$ws.Range("R418").Value = "Hortaliza"
